$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove every horizontal-rule paragraph.
#    Each one is a standalone <w:p> whose only content is a <w:pict><v:rect
#    .../></w:pict> "hr" drawing (o:hr="t"). Detect them via WordOpenXML and
#    delete the whole paragraph (Range.Delete also removes the paragraph
#    mark, so no empty paragraph is left behind). Walk back-to-front so
#    earlier indices stay valid while later ones are removed.
# ---------------------------------------------------------------------------
$hrParagraphs = @()
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.WordOpenXML -like "*o:hr=*") {
        $hrParagraphs += $i
    }
}

for ($j = $hrParagraphs.Count - 1; $j -ge 0; $j--) {
    $idx = $hrParagraphs[$j]
    $d.Paragraphs.Item($idx).Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) The paragraph that used to follow the very last horizontal rule
#    ("Last Updated: December 17, 2025", now the final paragraph in the
#    document) switches from the "FirstParagraph" style to "BodyText".
#    Changing the style clears direct character formatting, so re-apply the
#    italics the run already carried.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Style = "Body Text"
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastText = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastText.Font.Italic = $true

# ---------------------------------------------------------------------------
# 3) Fix the mangled multiplication sign: the two runs that used "×" now use
#    the mis-decoded "Ã—" sequence instead.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("×", $false, $false, $false, $false, $false, $true, 1, $false, "Ã—", 2) | Out-Null
